$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tr = $s.Shapes.Item(2).TextFrame.TextRange
$tr.InsertBefore("pigigpgpp")
$full = $s.Shapes.Item(2).TextFrame.TextRange
$newRun = $full.Characters(1, 9)
$newRun.LanguageID = "pt-BR"
